# Apply the cryptos-list price/volume refresh described by the commit diff.
# Column D (Price) values are forced to Text via NumberFormat '@' before the
# assignment so number-like strings (e.g. '1.001', '0.9999', '1.0000') are not
# auto-coerced to numeric values by Excel and keep their exact literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.011.01'
$ws.Range('E2').Value = '  -0.03%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.885.90'
$ws.Range('E3').Value = '  -1.60%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.20%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.54'
$ws.Range('E5').Value = '  -2.42%  '

# Row 6
$ws.Range('E6').Value = '  +0.17%  '

# Row 7
$ws.Range('E7').Value = '  -3.27%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4055'
$ws.Range('E8').Value = '  -0.20%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.69'
$ws.Range('E9').Value = '  -0.68%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07965'
$ws.Range('E10').Value = '  -2.83%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9913'
$ws.Range('E11').Value = '  -4.22%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.65'
$ws.Range('E12').Value = '  -4.14%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.884.72'
$ws.Range('E13').Value = '  -0.98%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.905'
$ws.Range('E14').Value = '  -3.42%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.074'
$ws.Range('E15').Value = '  -4.51%  '

# Row 16
$ws.Range('E16').Value = '  +0.15%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.35'
$ws.Range('E17').Value = '  -3.65%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001029'
$ws.Range('E18').Value = '  -2.50%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06549'
$ws.Range('E19').Value = '  -1.08%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.42'
$ws.Range('E20').Value = '  -3.07%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  +0.06%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.030.64'
$ws.Range('E22').Value = '  -0.02%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.423'
$ws.Range('E23').Value = '  -2.80%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.40'
$ws.Range('E24').Value = '  +1.51%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.205'
$ws.Range('E25').Value = '  -2.78%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.119.16'
$ws.Range('E26').Value = '  -0.51%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.62'
$ws.Range('E27').Value = '  -2.62%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.57'
$ws.Range('E28').Value = '  -2.45%  '

# Row 29
$ws.Range('E29').Value = '  -4.19%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.447'
$ws.Range('E30').Value = '  -2.25%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.70'
$ws.Range('E31').Value = '  -2.88%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.001'
$ws.Range('E32').Value = '  -2.22%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09325'
$ws.Range('E33').Value = '  -2.80%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.602'
$ws.Range('E34').Value = '  -1.45%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.409'
$ws.Range('E35').Value = '  -1.91%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.277'
$ws.Range('E36').Value = '  -3.04%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06055'
$ws.Range('E37').Value = '  -2.80%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02222'
$ws.Range('E38').Value = '  -3.06%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.275'
$ws.Range('E39').Value = '  -5.01%  '

# Row 40
$ws.Range('E40').Value = '  -2.74%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9999'
$ws.Range('E41').Value = '  +0.16%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5779'
$ws.Range('E42').Value = '  -4.76%  '

# Row 43
$ws.Range('E43').Value = '  -4.46%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.11'
$ws.Range('E44').Value = '  -5.04%  '

# Row 45
$ws.Range('E45').Value = '  -1.85%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07482'
$ws.Range('E46').Value = '  +2.44%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.04'
$ws.Range('E47').Value = '  -2.78%  '

# Row 48
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.257'
$ws.Range('E48').Value = '  +4.09%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5449'
$ws.Range('E49').Value = '  -3.77%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.898'
$ws.Range('E50').Value = '  -4.60%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.08'
$ws.Range('E51').Value = '  -1.99%  '
